$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "batch_size_1.0"

# Update row 2 values
$ws.Range("B2").Value = 0.6521704196929932
$ws.Range("E2").Value = 0.1123666666666667
$ws.Range("F2").Value = 0.1135

# Update row 3 values
$ws.Range("B3").Value = 0.426203727722168
$ws.Range("C3").Value = 20.27933154595121
$ws.Range("D3").Value = 21.38379761915692
$ws.Range("E3").Value = 0.9769833333333333
$ws.Range("F3").Value = 0.9805
